$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 = 14, Q1 = 15, matching formatting of O1 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null

# --- Data rows (rows 2-25): same new pattern for I:Q on every row ---
# I=2, J=2, K=1, L=2, M=2, N=2, O=1, P=2, Q=2
$rowCount = 24
$colCount = 9
$data = New-Object 'object[,]' $rowCount,$colCount
$rowPattern = @(2,2,1,2,2,2,1,2,2)
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $data[$r,$c] = $rowPattern[$c]
    }
}
$ws.Range("I2:Q25").Value = $data

$excel.CutCopyMode = 0
